# "Indirect System" feature tracker update.
#
# Two new planned features are introduced ("Формат задания сцены" and
# "Плагин для Blender"), two previously-planned features ("Physics" and
# "Screenshots") move from the "planned" (C) column into the
# "implemented" (B) column, and the "planned" (C) column is refilled with
# the features that take their place (including the two brand-new ones).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C ("Планируется реализовать" / planned) gets new contents ---
# Row2: Physics -> Shading (using shadow map)
$ws.Range("C2").Value = "Shading (using shadow map)"
# Row3: Screenshots -> Audio support
$ws.Range("C3").Value = "Audio support"
# Row4: Shading (using shadow map) -> Particle system
$ws.Range("C4").Value = "Particle system"
# Row5: Audio support -> Формат задания сцены (new)
$ws.Range("C5").Value = "Формат задания сцены"
# Row6: Particle system -> Плагин для Blender (new)
$ws.Range("C6").Value = "Плагин для Blender"

# --- Physics & Screenshots move into column B ("Реализовано" / implemented) ---
$ws.Range("B10").Value = "Physics"
$ws.Range("B10").Interior.Color = $ws.Range("B9").Interior.Color

$ws.Range("B11").Value = "Screenshots"
$ws.Range("B11").Interior.Color = $ws.Range("B9").Interior.Color

# --- The matching rows in column A ("Все фичи") flip from the planned ---
# --- (yellow) fill to the implemented (green) fill ---
$ws.Range("A6").Interior.Color = $ws.Range("B6").Interior.Color
$ws.Range("A9").Interior.Color = $ws.Range("B9").Interior.Color

# --- Two new rows appended to column A for the brand-new planned features ---
$ws.Range("A26").Value = "Формат задания сцены"
$ws.Range("A26").Interior.Color = $ws.Range("A20").Interior.Color

$ws.Range("A27").Value = "Плагин для Blender"
$ws.Range("A27").Interior.Color = $ws.Range("A20").Interior.Color

# --- Selection moves to C7 ---
$ws.Range("C7").Select()
